$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3874
$ws.Range("E2").Value = 317
$ws.Range("F2").Value = 317
$ws.Range("G2").Value = 276
$ws.Range("H2").Value = 256
$ws.Range("I2").Value = 256
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 10345
$ws.Range("L2").Value = 2943
$ws.Range("M2").Value = 7402
$ws.Range("N2").Value = 7390
$ws.Range("O2").Value = 12
$ws.Range("P2").Value = 232
$ws.Range("Q2").Value = 348
$ws.Range("R2").Value = -309
$ws.Range("S2").Value = -81
$ws.Range("T2").Value = 274
$ws.Range("U2").Value = 74
$ws.Range("V2").Value = 1260
$ws.Range("W2").Value = 8.18
$ws.Range("X2").Value = 6.62
$ws.Range("Y2").Value = 3.5
$ws.Range("Z2").Value = 2.49
$ws.Range("AA2").Value = 39.76
$ws.Range("AB2").Value = 3079.97
$ws.Range("AC2").Value = 5530
$ws.Range("AD2").Value = 18.81
$ws.Range("AE2").Value = 159349
$ws.Range("AF2").Value = 0.65
$ws.Range("AG2").Value = 2000
$ws.Range("AH2").Value = 1.92
$ws.Range("AI2").Value = 36.16
$ws.Range("AJ2").Value = 4637790

# Row 3
$ws.Range("D3").Value = 4491
$ws.Range("E3").Value = 557
$ws.Range("F3").Value = 557
$ws.Range("G3").Value = 514
$ws.Range("H3").Value = 476
$ws.Range("I3").Value = 476
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 10847
$ws.Range("L3").Value = 3052
$ws.Range("M3").Value = 7795
$ws.Range("N3").Value = 7783
$ws.Range("O3").Value = 12
$ws.Range("P3").Value = 232
$ws.Range("Q3").Value = 48
$ws.Range("R3").Value = -61
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 53
$ws.Range("U3").Value = -5
$ws.Range("V3").Value = 1356
$ws.Range("W3").Value = 12.4
$ws.Range("X3").Value = 10.6
$ws.Range("Y3").Value = 6.27
$ws.Range("Z3").Value = 4.49
$ws.Range("AA3").Value = 39.15
$ws.Range("AB3").Value = 3237.52
$ws.Range("AC3").Value = 10264
$ws.Range("AD3").Value = 15
$ws.Range("AE3").Value = 167826
$ws.Range("AF3").Value = 0.92
$ws.Range("AG3").Value = 2000
$ws.Range("AH3").Value = 1.3
$ws.Range("AI3").Value = 19.49
$ws.Range("AJ3").Value = 4637790

# Row 4
$ws.Range("D4").Value = 4978
$ws.Range("E4").Value = 869
$ws.Range("F4").Value = 869
$ws.Range("G4").Value = 843
$ws.Range("H4").Value = 794
$ws.Range("I4").Value = 794
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 11719
$ws.Range("L4").Value = 3304
$ws.Range("M4").Value = 8415
$ws.Range("N4").Value = 8402
$ws.Range("O4").Value = 12
$ws.Range("P4").Value = 232
$ws.Range("Q4").Value = 171
$ws.Range("R4").Value = -125
$ws.Range("S4").Value = -1
$ws.Range("T4").Value = 53
$ws.Range("U4").Value = 118
$ws.Range("V4").Value = 1449
$ws.Range("W4").Value = 17.45
$ws.Range("X4").Value = 15.95
$ws.Range("Y4").Value = 9.81
$ws.Range("Z4").Value = 7.04
$ws.Range("AA4").Value = 39.27
$ws.Range("AB4").Value = 3521.41
$ws.Range("AC4").Value = 17111
$ws.Range("AD4").Value = 7.04
$ws.Range("AE4").Value = 181174
$ws.Range("AF4").Value = 0.67
$ws.Range("AG4").Value = 2000
$ws.Range("AH4").Value = 1.66
$ws.Range("AI4").Value = 11.69
$ws.Range("AJ4").Value = 4637790

# Row 5
$ws.Range("D5").Value = 4966
$ws.Range("E5").Value = 471
$ws.Range("F5").Value = 471
$ws.Range("G5").Value = 474
$ws.Range("H5").Value = 407
$ws.Range("I5").Value = 406
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 11453
$ws.Range("L5").Value = 2821
$ws.Range("M5").Value = 8633
$ws.Range("N5").Value = 8620
$ws.Range("O5").Value = 13
$ws.Range("P5").Value = 232
$ws.Range("Q5").Value = 588
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = -606
$ws.Range("T5").Value = 68
$ws.Range("U5").Value = 520
$ws.Range("V5").Value = 937
$ws.Range("W5").Value = 9.48
$ws.Range("X5").Value = 8.2
$ws.Range("Y5").Value = 4.78
$ws.Range("Z5").Value = 3.51
$ws.Range("AA5").Value = 32.67
$ws.Range("AB5").Value = 3654.99
$ws.Range("AC5").Value = 8764
$ws.Range("AD5").Value = 13.06
$ws.Range("AE5").Value = 185858
$ws.Range("AF5").Value = 0.62
$ws.Range("AG5").Value = 2000
$ws.Range("AH5").Value = 1.75
$ws.Range("AI5").Value = 22.82
$ws.Range("AJ5").Value = 4637790

# Row 6
$ws.Range("D6").Value = 5432
$ws.Range("E6").Value = 411
$ws.Range("F6").Value = 411
$ws.Range("G6").Value = 390
$ws.Range("H6").Value = 468
$ws.Range("I6").Value = 467
$ws.Range("K6").Value = 11948
$ws.Range("L6").Value = 2989
$ws.Range("M6").Value = 8959
$ws.Range("N6").Value = 8946
$ws.Range("P6").Value = 232
$ws.Range("Q6").Value = 156
$ws.Range("R6").Value = -101
$ws.Range("S6").Value = -34
$ws.Range("T6").Value = 50
$ws.Range("U6").Value = 105
$ws.Range("V6").Value = 996
$ws.Range("W6").Value = 7.57
$ws.Range("X6").Value = 8.61
$ws.Range("Y6").Value = 5.32
$ws.Range("Z6").Value = 4
$ws.Range("AA6").Value = 33.37
$ws.Range("AB6").Value = 3793.05
$ws.Range("AC6").Value = 10077
$ws.Range("AD6").Value = 6.94
$ws.Range("AE6").Value = 192883
$ws.Range("AF6").Value = 0.36
$ws.Range("AG6").Value = 2000
$ws.Range("AH6").Value = 2.86
$ws.Range("AI6").Value = 19.85
$ws.Range("AJ6").Value = 4637790

# Remove forecast-year rows (7-9): drop all D:AJ cell content, keep A/B/C labels
$ws.Range("D7:AJ9").ClearContents()

